# Scheduled data refresh: update market-price / leve-profit figures across sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1952.4445
$ws.Cells.Item(17, 10).Value = 2749.4546
$ws.Cells.Item(17, 12).Value = 8248.363799999999
$ws.Cells.Item(17, 14).Value = -8584.363799999999
$ws.Cells.Item(32, 8).Value = 977.5
$ws.Cells.Item(32, 10).Value = 970.3333
$ws.Cells.Item(32, 12).Value = 970.3333
$ws.Cells.Item(32, 14).Value = -1622.3333
$ws.Cells.Item(33, 8).Value = 109.125
$ws.Cells.Item(33, 9).Value = 83.09090999999999
$ws.Cells.Item(33, 11).Value = 83.09090999999999
$ws.Cells.Item(33, 13).Value = 145.90909
$ws.Cells.Item(40, 8).Value = 6019.727
$ws.Cells.Item(40, 9).Value = 4681.4287
$ws.Cells.Item(40, 10).Value = 8361.75
$ws.Cells.Item(40, 11).Value = 4681.4287
$ws.Cells.Item(40, 12).Value = 8361.75
$ws.Cells.Item(40, 13).Value = -4506.4287
$ws.Cells.Item(40, 14).Value = -8711.75
$ws.Cells.Item(59, 8).Value = 58.5
$ws.Cells.Item(59, 9).Value = 58.5
$ws.Cells.Item(59, 11).Value = 175.5
$ws.Cells.Item(59, 13).Value = 381.5
$ws.Cells.Item(129, 8).Value = 694.86664
$ws.Cells.Item(129, 9).Value = 530.2143
$ws.Cells.Item(129, 10).Value = 3000
$ws.Cells.Item(129, 11).Value = 1590.6429
$ws.Cells.Item(129, 12).Value = 9000
$ws.Cells.Item(129, 13).Value = 3409.3571
$ws.Cells.Item(129, 14).Value = -19000

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value = 2904.7144
$ws.Cells.Item(38, 9).Value = 2555.6667
$ws.Cells.Item(38, 10).Value = 4999
$ws.Cells.Item(38, 11).Value = 2555.6667
$ws.Cells.Item(38, 12).Value = 4999
$ws.Cells.Item(38, 13).Value = -2088.6667
$ws.Cells.Item(38, 14).Value = -5933
$ws.Cells.Item(61, 8).Value = 2043.9286
$ws.Cells.Item(61, 9).Value = 1931.9231
$ws.Cells.Item(61, 10).Value = 3500
$ws.Cells.Item(61, 11).Value = 1931.9231
$ws.Cells.Item(61, 12).Value = 3500
$ws.Cells.Item(61, 13).Value = -1719.9231
$ws.Cells.Item(61, 14).Value = -3924
$ws.Cells.Item(63, 8).Value = 2157.5
$ws.Cells.Item(63, 9).Value = 2157.5
$ws.Cells.Item(63, 11).Value = 2157.5
$ws.Cells.Item(63, 13).Value = -1471.5
$ws.Cells.Item(66, 8).Value = 2157.5
$ws.Cells.Item(66, 9).Value = 2157.5
$ws.Cells.Item(66, 11).Value = 10787.5
$ws.Cells.Item(66, 13).Value = -7355.5
$ws.Cells.Item(74, 8).Value = 4047.4119
$ws.Cells.Item(74, 9).Value = 3965.0715
$ws.Cells.Item(74, 10).Value = 4431.6665
$ws.Cells.Item(74, 11).Value = 3965.0715
$ws.Cells.Item(74, 12).Value = 4431.6665
$ws.Cells.Item(74, 13).Value = -3091.0715
$ws.Cells.Item(74, 14).Value = -6179.6665
$ws.Cells.Item(77, 8).Value = 4047.4119
$ws.Cells.Item(77, 9).Value = 3965.0715
$ws.Cells.Item(77, 10).Value = 4431.6665
$ws.Cells.Item(77, 11).Value = 19825.3575
$ws.Cells.Item(77, 12).Value = 22158.3325
$ws.Cells.Item(77, 13).Value = -15457.3575
$ws.Cells.Item(77, 14).Value = -30894.3325
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 2043.9286
$ws.Cells.Item(136, 9).Value = 1931.9231
$ws.Cells.Item(136, 10).Value = 3500
$ws.Cells.Item(136, 11).Value = 5795.7693
$ws.Cells.Item(136, 12).Value = 10500
$ws.Cells.Item(136, 13).Value = -3245.7693
$ws.Cells.Item(136, 14).Value = -15600

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 650
$ws.Cells.Item(11, 9).Value = 650
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 650
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -510
$ws.Cells.Item(11, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 750
$ws.Cells.Item(20, 9).Value = 500
$ws.Cells.Item(20, 11).Value = 500
$ws.Cells.Item(20, 13).Value = -253
$ws.Cells.Item(44, 8).Value = 13275
$ws.Cells.Item(44, 10).Value = 13275
$ws.Cells.Item(44, 12).Value = 13275
$ws.Cells.Item(44, 14).Value = -14269
$ws.Cells.Item(86, 8).Value = 2719.2
$ws.Cells.Item(86, 9).Value = 1358.75
$ws.Cells.Item(86, 11).Value = 1358.75
$ws.Cells.Item(86, 13).Value = -235.75
$ws.Cells.Item(89, 8).Value = 2719.2
$ws.Cells.Item(89, 9).Value = 1358.75
$ws.Cells.Item(89, 11).Value = 6793.75
$ws.Cells.Item(89, 13).Value = -1177.75
$ws.Cells.Item(111, 8).Value = 40000
$ws.Cells.Item(111, 10).Value = 40000
$ws.Cells.Item(111, 12).Value = 40000
$ws.Cells.Item(111, 14).Value = -48180

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 10000
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(7, 8).Value = 72.9375
$ws.Cells.Item(7, 9).Value = 53.75
$ws.Cells.Item(7, 10).Value = 130.5
$ws.Cells.Item(7, 11).Value = 53.75
$ws.Cells.Item(7, 12).Value = 130.5
$ws.Cells.Item(7, 13).Value = 59.25
$ws.Cells.Item(7, 14).Value = -356.5
$ws.Cells.Item(31, 8).Value = 3239.325
$ws.Cells.Item(31, 9).Value = 1878.65
$ws.Cells.Item(31, 11).Value = 1878.65
$ws.Cells.Item(31, 13).Value = -1583.65
$ws.Cells.Item(33, 8).Value = 4059.8
$ws.Cells.Item(33, 9).Value = 1949.75
$ws.Cells.Item(33, 10).Value = 12500
$ws.Cells.Item(33, 11).Value = 1949.75
$ws.Cells.Item(33, 12).Value = 12500
$ws.Cells.Item(33, 13).Value = -1570.75
$ws.Cells.Item(33, 14).Value = -13258
$ws.Cells.Item(34, 8).Value = 3239.325
$ws.Cells.Item(34, 9).Value = 1878.65
$ws.Cells.Item(34, 11).Value = 1878.65
$ws.Cells.Item(34, 13).Value = -1676.65
$ws.Cells.Item(62, 8).Value = 1833
$ws.Cells.Item(62, 10).Value = 1999.5
$ws.Cells.Item(62, 12).Value = 1999.5
$ws.Cells.Item(62, 14).Value = -3247.5
$ws.Cells.Item(65, 8).Value = 1833
$ws.Cells.Item(65, 10).Value = 1999.5
$ws.Cells.Item(65, 12).Value = 9997.5
$ws.Cells.Item(65, 14).Value = -16237.5
$ws.Cells.Item(132, 8).Value = 2377.7727
$ws.Cells.Item(132, 9).Value = 2348.1428
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 7044.428400000001
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -4514.428400000001
$ws.Cells.Item(132, 14).Value = -14060
$ws.Cells.Item(134, 8).Value = 2361.1
$ws.Cells.Item(134, 9).Value = 1301.9286
$ws.Cells.Item(134, 11).Value = 3905.7858
$ws.Cells.Item(134, 13).Value = -1370.7858
$ws.Cells.Item(138, 8).Value = 127996
$ws.Cells.Item(138, 10).Value = 127996
$ws.Cells.Item(138, 12).Value = 127996
$ws.Cells.Item(138, 14).Value = -138276

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 436.33334
$ws.Cells.Item(47, 9).Value = 423.75
$ws.Cells.Item(47, 10).Value = 461.5
$ws.Cells.Item(47, 11).Value = 1271.25
$ws.Cells.Item(47, 12).Value = 1384.5
$ws.Cells.Item(47, 13).Value = -840.25
$ws.Cells.Item(47, 14).Value = -2246.5
$ws.Cells.Item(58, 8).Value = 1751
$ws.Cells.Item(58, 10).Value = 1999.5
$ws.Cells.Item(58, 12).Value = 5998.5
$ws.Cells.Item(58, 14).Value = -6254.5
$ws.Cells.Item(129, 8).Value = 1957.3
$ws.Cells.Item(129, 9).Value = 778
$ws.Cells.Item(129, 11).Value = 2334
$ws.Cells.Item(129, 13).Value = 2666
$ws.Cells.Item(131, 8).Value = 2749.5
$ws.Cells.Item(131, 10).Value = 2999.3333
$ws.Cells.Item(131, 12).Value = 8997.999899999999
$ws.Cells.Item(131, 14).Value = -19077.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 71.23529000000001
$ws.Cells.Item(2, 9).Value = 60.18182
$ws.Cells.Item(2, 10).Value = 91.5
$ws.Cells.Item(2, 11).Value = 60.18182
$ws.Cells.Item(2, 12).Value = 91.5
$ws.Cells.Item(2, 13).Value = 52.81818
$ws.Cells.Item(2, 14).Value = -317.5
$ws.Cells.Item(80, 8).Value = 4825
$ws.Cells.Item(80, 9).Value = 3933.3333
$ws.Cells.Item(80, 11).Value = 3933.3333
$ws.Cells.Item(80, 13).Value = -2935.3333
$ws.Cells.Item(83, 8).Value = 4825
$ws.Cells.Item(83, 9).Value = 3933.3333
$ws.Cells.Item(83, 11).Value = 19666.6665
$ws.Cells.Item(83, 13).Value = -14674.6665

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2074.8333
$ws.Cells.Item(22, 9).Value = 1579.8
$ws.Cells.Item(22, 10).Value = 2428.4285
$ws.Cells.Item(22, 11).Value = 1579.8
$ws.Cells.Item(22, 12).Value = 2428.4285
$ws.Cells.Item(22, 13).Value = -1284.8
$ws.Cells.Item(22, 14).Value = -3018.4285
$ws.Cells.Item(27, 8).Value = 2074.8333
$ws.Cells.Item(27, 9).Value = 1579.8
$ws.Cells.Item(27, 10).Value = 2428.4285
$ws.Cells.Item(27, 11).Value = 1579.8
$ws.Cells.Item(27, 12).Value = 2428.4285
$ws.Cells.Item(27, 13).Value = -1472.8
$ws.Cells.Item(27, 14).Value = -2642.4285
$ws.Cells.Item(46, 8).Value = 2102.4
$ws.Cells.Item(46, 10).Value = 3382.182
$ws.Cells.Item(46, 12).Value = 3382.182
$ws.Cells.Item(46, 14).Value = -3758.182
$ws.Cells.Item(59, 8).Value = 33000
$ws.Cells.Item(59, 10).Value = 33000
$ws.Cells.Item(59, 12).Value = 33000
$ws.Cells.Item(59, 14).Value = -34308
$ws.Cells.Item(61, 8).Value = 71432536
$ws.Cells.Item(61, 9).Value = 142858350
$ws.Cells.Item(61, 11).Value = 142858350
$ws.Cells.Item(61, 13).Value = -142858148
$ws.Cells.Item(113, 8).Value = 71432536
$ws.Cells.Item(113, 9).Value = 142858350
$ws.Cells.Item(113, 11).Value = 142858350
$ws.Cells.Item(113, 13).Value = -142856180
$ws.Cells.Item(136, 8).Value = 1950.5
$ws.Cells.Item(136, 9).Value = 1754.3846
$ws.Cells.Item(136, 10).Value = 4500
$ws.Cells.Item(136, 11).Value = 5263.1538
$ws.Cells.Item(136, 12).Value = 13500
$ws.Cells.Item(136, 13).Value = -2713.1538
$ws.Cells.Item(136, 14).Value = -18600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 34333.75
$ws.Cells.Item(4, 9).Value = 37375.453
$ws.Cells.Item(4, 10).Value = 875
$ws.Cells.Item(4, 11).Value = 37375.453
$ws.Cells.Item(4, 12).Value = 875
$ws.Cells.Item(4, 13).Value = -37262.453
$ws.Cells.Item(4, 14).Value = -1101
$ws.Cells.Item(122, 8).Value = 2554.1667
$ws.Cells.Item(122, 9).Value = 1475.875
$ws.Cells.Item(122, 11).Value = 4427.625
$ws.Cells.Item(122, 13).Value = -1977.625
$ws.Cells.Item(126, 8).Value = 4007.35
$ws.Cells.Item(126, 9).Value = 1876.6923
$ws.Cells.Item(126, 11).Value = 5630.0769
$ws.Cells.Item(126, 13).Value = -3160.0769
$ws.Cells.Item(132, 8).Value = 2181.182
$ws.Cells.Item(132, 9).Value = 1924.25
$ws.Cells.Item(132, 10).Value = 2866.3333
$ws.Cells.Item(132, 11).Value = 5772.75
$ws.Cells.Item(132, 12).Value = 8598.999899999999
$ws.Cells.Item(132, 13).Value = -3242.75
$ws.Cells.Item(132, 14).Value = -13658.9999
$ws.Cells.Item(136, 8).Value = 2052.5186
$ws.Cells.Item(136, 9).Value = 1726.4
$ws.Cells.Item(136, 10).Value = 2984.2856
$ws.Cells.Item(136, 11).Value = 5179.200000000001
$ws.Cells.Item(136, 12).Value = 8952.856800000001
$ws.Cells.Item(136, 13).Value = -2629.200000000001
$ws.Cells.Item(136, 14).Value = -14052.8568

